# Updated cryptos list with GitHub Actions
# Applies the cell-value changes described by the diff against cryptos.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell while keeping it as
# plain text (no thousands-grouping / locale re-interpretation of values
# such as '1.00' or '606.65') and without leaving a lingering custom
# number format or quote-prefix style behind on the cell.
function Set-TextValue($cell, $text) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = '@'
    $rng.Value = $text
    $rng.Style = 'Normal'
}

Set-TextValue 'D2' '65.713.03'
Set-TextValue 'E2' '  +3.00%  '
Set-TextValue 'D3' '2.667.02'
Set-TextValue 'E3' '  +1.66%  '
Set-TextValue 'E4' '  +0.00%  '
Set-TextValue 'D5' '606.65'
Set-TextValue 'E5' '  +1.85%  '
Set-TextValue 'D6' '158.56'
Set-TextValue 'E6' '  +4.87%  '
Set-TextValue 'E7' '  -0.04%  '
Set-TextValue 'E8' '  +0.75%  '
Set-TextValue 'E9' '  +8.23%  '
Set-TextValue 'E10' '  +2.38%  '
Set-TextValue 'D11' '5.85'
Set-TextValue 'E11' '  +0.63%  '
Set-TextValue 'E12' '  +1.76%  '
Set-TextValue 'D13' '29.85'
Set-TextValue 'E13' '  +7.14%  '
Set-TextValue 'D14' '0.0000194'
Set-TextValue 'E14' '  +14.22%  '
Set-TextValue 'D15' '3.148.30'
Set-TextValue 'E15' '  +1.74%  '
Set-TextValue 'D16' '65.433.83'
Set-TextValue 'E16' '  +2.68%  '
Set-TextValue 'D17' '2.668.39'
Set-TextValue 'E17' '  +1.27%  '
Set-TextValue 'E18' '  +5.00%  '
Set-TextValue 'E19' '  +2.58%  '
Set-TextValue 'D20' '360.46'
Set-TextValue 'E20' '  +3.80%  '
Set-TextValue 'D21' '7.39'
Set-TextValue 'E21' '  +5.38%  '
Set-TextValue 'E22' '  +0.06%  '
Set-TextValue 'D23' '69.07'
Set-TextValue 'E23' '  +2.65%  '
Set-TextValue 'D24' '1.70'
Set-TextValue 'E24' '  +1.26%  '
Set-TextValue 'D25' '9.59'
Set-TextValue 'E25' '  +4.86%  '
Set-TextValue 'E26' '  +17.50%  '
Set-TextValue 'E27' '  -1.61%  '
Set-TextValue 'D28' '8.24'
Set-TextValue 'E28' '  -0.54%  '
Set-TextValue 'E29' '  +2.03%  '
Set-TextValue 'B30' 'Binance-PegBSC-USD'
Set-TextValue 'C30' 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue 'D30' '1.00'
Set-TextValue 'E30' '  +0.25%  '
Set-TextValue 'B31' 'PancakeSwap'
Set-TextValue 'C31' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D31' '2.19'
Set-TextValue 'E31' '  +5.99%  '
Set-TextValue 'D32' '541.01'
Set-TextValue 'E32' '  -1.87%  '
Set-TextValue 'D33' '1.85'
Set-TextValue 'E33' '  +2.12%  '
Set-TextValue 'E34' '  +5.24%  '
Set-TextValue 'D35' '6.38'
Set-TextValue 'E35' '  +4.50%  '
Set-TextValue 'E36' '  +4.26%  '
Set-TextValue 'E37' '  +3.70%  '
Set-TextValue 'D38' '2.03'
Set-TextValue 'E38' '  +2.76%  '
Set-TextValue 'D39' '162.92'
Set-TextValue 'E39' '  -0.52%  '
Set-TextValue 'D40' '1.00'
Set-TextValue 'E40' '  +0.06%  '
Set-TextValue 'E41' '  +0.01%  '
Set-TextValue 'D42' '42.49'
Set-TextValue 'E42' '  +6.55%  '
Set-TextValue 'D43' '167.09'
Set-TextValue 'E43' '  -0.42%  '
Set-TextValue 'E44' '  +2.22%  '
Set-TextValue 'D45' '2.38'
Set-TextValue 'E45' '  +8.62%  '
Set-TextValue 'D46' '0.0616'
Set-TextValue 'E46' '  +5.61%  '
Set-TextValue 'D47' '23.25'
Set-TextValue 'E47' '  -0.21%  '
Set-TextValue 'D48' '0.664'
Set-TextValue 'E48' '  +4.11%  '
Set-TextValue 'E49' '  +5.75%  '
Set-TextValue 'E50' '  +2.47%  '
Set-TextValue 'D51' '19.94'
Set-TextValue 'E51' '  +3.60%  '
